$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: drop the stray "no-fill" explicit formatting that was left on the
#     separator / StandUp cells (G3:X3). The separator cells (G,I,K,...) never
#     held content, so they disappear entirely; the cells that still hold text
#     ("StandUp x3", "D-DAY") just lose their pointless explicit format.
foreach ($col in @("G","I","K","M","O","Q","S","U","W")) {
    $ws.Range($col + "3").Clear()
}
foreach ($col in @("H","J","L","N","P","R","T","V","X")) {
    $ws.Range($col + "3").ClearFormats()
}

# --- Row 4: clear out the old blanket formatting across the row, then lay
#     down the new content - "Search in dummy DB" / "Integrate DB" tasks that
#     now live under Day 2 / Day 3, with an orange marker cell before them.
$ws.Range("G4:X4").Clear()
$ws.Range("G4").Interior.Color = 49407
$ws.Range("H4").Value = "Search in dummy DB"
$ws.Range("J4").Value = "Integrate DB"

# --- Row 5: same treatment - orange marker then the "Find API" task.
$ws.Range("G5:X5").Clear()
$ws.Range("G5").Interior.Color = 49407
$ws.Range("H5").Value = "Find API"

# --- Row 6: brand new orange marker + "Find Database" task cell.
$ws.Range("G6").Interior.Color = 49407
$ws.Range("H6").Value = "Find Database"

# --- Row 14: "Create a dummy DB" marker switches from orange to green.
$ws.Range("E14").Interior.Color = 5287936

# --- Reflect where the user left the selection after these edits.
$ws.Range("O4").Select() | Out-Null

Write-Output "edit applied"
